$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the header labels in row 1 (B1:F1) per the commit's finalized
# input list ordering. A1 (bedrooms_1) stays the same.
$ws.Range("B1").Value = "living_rooms_1"
$ws.Range("C1").Value = "bedrooms_2"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "kitchens_1"
$ws.Range("F1").Value = "kitchens_2"
